# Apply updated cryptocurrency price/volume data (and two row re-labels)
# Each target cell is forced to Text before assignment (NumberFormat "@"),
# then ClearFormats() restores the default (unstyled) cell formatting so only
# the cell VALUE changes and no stray styling is introduced - this matches how
# the source file stores these as literal text (e.g. "30.277.80", "  -0.22%  ")
# rather than as numbers, even though some look numeric.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue 'D2' '30.301.83'
Set-TextValue 'E2' '  -0.20%  '
Set-TextValue 'D3' '1.932.68'
Set-TextValue 'E3' '  -0.49%  '
Set-TextValue 'D4' '0.9997'
Set-TextValue 'E4' '  -0.14%  '
Set-TextValue 'D5' '0.7504'
Set-TextValue 'E5' '  +4.11%  '
Set-TextValue 'D6' '250.11'
Set-TextValue 'E6' '  -0.53%  '
Set-TextValue 'D7' '0.9987'
Set-TextValue 'E7' '  -0.20%  '
Set-TextValue 'D8' '0.3234'
Set-TextValue 'E8' '  -3.13%  '
Set-TextValue 'D9' '28.13'
Set-TextValue 'E9' '  -2.71%  '
Set-TextValue 'D10' '0.07149'
Set-TextValue 'E10' '  -3.71%  '
Set-TextValue 'D11' '0.7923'
Set-TextValue 'E11' '  -2.85%  '
Set-TextValue 'D12' '0.08005'
Set-TextValue 'E12' '  -1.55%  '
Set-TextValue 'D13' '1.933.13'
Set-TextValue 'E13' '  -0.41%  '
Set-TextValue 'D14' '5.404'
Set-TextValue 'E14' '  -1.53%  '
Set-TextValue 'D15' '94.60'
Set-TextValue 'E15' '  -0.37%  '
Set-TextValue 'D16' '14.58'
Set-TextValue 'E16' '  -2.87%  '
Set-TextValue 'D17' '30.305.57'
Set-TextValue 'E17' '  -0.22%  '
Set-TextValue 'D18' '253.98'
Set-TextValue 'E18' '  +0.41%  '
Set-TextValue 'D19' '0.000008060'
Set-TextValue 'E19' '  -3.05%  '
Set-TextValue 'D20' '5.785'
Set-TextValue 'E20' '  -1.84%  '
Set-TextValue 'D21' '2.187.93'
Set-TextValue 'E21' '  -0.27%  '
Set-TextValue 'D22' '0.9991'
Set-TextValue 'E22' '  -0.14%  '
Set-TextValue 'D23' '0.9994'
Set-TextValue 'E23' '  -0.23%  '
Set-TextValue 'E24' '  -2.05%  '
Set-TextValue 'D25' '9.598'
Set-TextValue 'E25' '  -2.60%  '
Set-TextValue 'D26' '164.42'
Set-TextValue 'E26' '  +0.98%  '
Set-TextValue 'D27' '2.330'
Set-TextValue 'E27' '  -3.39%  '
Set-TextValue 'D28' '0.1346'
Set-TextValue 'E28' '  +2.32%  '
Set-TextValue 'D29' '19.13'
Set-TextValue 'E29' '  -1.56%  '
Set-TextValue 'D30' '1.360'
Set-TextValue 'E30' '  +0.93%  '
Set-TextValue 'E31' '  -2.65%  '
Set-TextValue 'D32' '4.435'
Set-TextValue 'E32' '  -0.64%  '
Set-TextValue 'D33' '4.155'
Set-TextValue 'E33' '  -1.93%  '
Set-TextValue 'D34' '1.302'
Set-TextValue 'E34' '  -0.77%  '
Set-TextValue 'D35' '0.05117'
Set-TextValue 'E35' '  -2.93%  '
Set-TextValue 'D36' '0.7499'
Set-TextValue 'E36' '  -0.72%  '
Set-TextValue 'D37' '2.761'
Set-TextValue 'E37' '  +0.56%  '
Set-TextValue 'E38' '  -1.26%  '
Set-TextValue 'D39' '2.800'
Set-TextValue 'E39' '  -1.79%  '
Set-TextValue 'D40' '78.37'
Set-TextValue 'E40' '  -4.76%  '
Set-TextValue 'D41' '6.411'
Set-TextValue 'E41' '  -2.91%  '
Set-TextValue 'D42' '0.4528'
Set-TextValue 'E42' '  -0.83%  '
Set-TextValue 'D43' '2.004'
Set-TextValue 'E43' '  -1.67%  '
Set-TextValue 'D44' '0.8427'
Set-TextValue 'E44' '  -0.63%  '
Set-TextValue 'D45' '0.9990'
Set-TextValue 'E45' '  -0.18%  '
Set-TextValue 'D46' '101.92'
Set-TextValue 'E46' '  -0.66%  '
Set-TextValue 'D47' '9.838'
Set-TextValue 'E47' '  +0.00%  '
Set-TextValue 'B48' 'Aptos'
Set-TextValue 'C48' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D48' '7.538'
Set-TextValue 'E48' '  +0.16%  '
Set-TextValue 'B49' 'Maker'
Set-TextValue 'C49' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D49' '993.41'
Set-TextValue 'E49' '  +12.49%  '
Set-TextValue 'D50' '37.48'
Set-TextValue 'E50' '  +1.31%  '
Set-TextValue 'B51' 'Decentraland'
Set-TextValue 'C51' 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue 'D51' '0.4178'
Set-TextValue 'E51' '  -0.58%  '
